$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, matching the style of the existing header cells (copy format from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells H2:H4, numeric value 0 (no special style, like F2:G4)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
